# Applies the "Add files via upload / 15/10" edit to the
# "Assignemnts-original" sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assignemnts-original")

# --- 1. Remove the "Function PaybyPaypal" entry (row 10, column A) ---
$ws.Range("A10").Value = $null

# --- 2. Fix the Full Name typo for roll number SE05407 (row 14) ---
$ws.Range("G14").Value = "Huỳnh Chí Cường"

# --- 3. Clear out the old "Function Cookies" row (row 19) ---
$ws.Range("A19:I19").Value = $null
$ws.Range("M19").Value = $null
$ws.Range("P19").Value = $null

# --- 4. Fill in the new assignment rows (20-25) ---
# Row 20: render to shoping
$ws.Range("A20").Value = "render to shoping"
$ws.Range("B20").Value = "GameStore"
$ws.Range("C20").Value = "simple"
$ws.Range("E20").Value = "Iteration 1"
$ws.Range("F20").Value = "SE150623"
$ws.Range("G20").Value = "Phan Thiên Ân"
$ws.Range("H20").Value = "T5"
$ws.Range("I20").Value = "Iteration 1"
$ws.Range("J20").Value = "simple"
$ws.Range("M20").Value = "simple"
$ws.Range("P20").Value = 60

# Row 21: dashboard
$ws.Range("A21").Value = "dashboard"
$ws.Range("B21").Value = "GameStore"
$ws.Range("C21").Value = "medium"
$ws.Range("E21").Value = "Iteration 1"
$ws.Range("F21").Value = "SE150679"
$ws.Range("G21").Value = "Huỳnh Chí Cường"
$ws.Range("H21").Value = "T5"
$ws.Range("I21").Value = "Iteration 1"
$ws.Range("J21").Value = "medium"
$ws.Range("M21").Value = "medium"
$ws.Range("P21").Value = 120

# Row 22: paypal
$ws.Range("A22").Value = "paypal"
$ws.Range("B22").Value = "GameStore"
$ws.Range("C22").Value = "complex"
$ws.Range("E22").Value = "Iteration 1"
$ws.Range("F22").Value = "SE150679"
$ws.Range("G22").Value = "Huỳnh Chí Cường"
$ws.Range("H22").Value = "T5"
$ws.Range("I22").Value = "Iteration 1"
$ws.Range("J22").Value = "complex"
$ws.Range("M22").Value = "complex"
$ws.Range("P22").Value = 180

# Row 23: contact
$ws.Range("A23").Value = "contact"
$ws.Range("B23").Value = "GameStore"
$ws.Range("C23").Value = "simple"
$ws.Range("E23").Value = "Iteration 1"
$ws.Range("F23").Value = "SE150674"
$ws.Range("G23").Value = "Võ Chí Cường"
$ws.Range("H23").Value = "T5"
$ws.Range("I23").Value = "Iteration 1"
$ws.Range("J23").Value = "simple"
$ws.Range("M23").Value = "simple"
$ws.Range("P23").Value = 60

# Row 24: search product by admin
$ws.Range("A24").Value = "search product by admin"
$ws.Range("B24").Value = "GameStore"
$ws.Range("C24").Value = "simple"
$ws.Range("E24").Value = "Iteration 1"
$ws.Range("F24").Value = "SE140461"
$ws.Range("G24").Value = "Nguyễn Quang Minh"
$ws.Range("H24").Value = "T5"
$ws.Range("I24").Value = "Iteration 1"
$ws.Range("J24").Value = "simple"
$ws.Range("M24").Value = "simple"
$ws.Range("P24").Value = 60

# Row 25: render contact
$ws.Range("A25").Value = "render contact"
$ws.Range("B25").Value = "GameStore"
$ws.Range("C25").Value = "simple"
$ws.Range("E25").Value = "Iteration 1"
$ws.Range("F25").Value = "SE63160"
$ws.Range("G25").Value = "Phạm Quang Quý"
$ws.Range("H25").Value = "T5"
$ws.Range("I25").Value = "Iteration 1"
$ws.Range("J25").Value = "simple"
$ws.Range("M25").Value = "simple"
$ws.Range("P25").Value = 60

# --- 5. Formatting touch-ups ---
# New rows use the same "Input" cell style as the rest of the table.
$ws.Range("A20:P25").Style = "Input"

# Column A needs to be a little wider to fit the new longer labels.
$ws.Columns.Item(1).ColumnWidth = 22.88671875

$wb.Save()
